# Rename the last header (H1) and add three new trailing headers (I1:K1)
# for product-image/meta SEO fields, matching the header row's existing
# look (same font/wrap formatting as the rest of row 1). Also add a blank
# second row (data row) below the headers with word-wrap enabled, matching
# how the refreshed template ships.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H1 used to read "Tên đường dẫn" -> now "Đường dẫn ảnh" (supports pulling
# product images from an outside server instead of just a local path).
$ws.Range("H1").Value = "Đường dẫn ảnh"

# New SEO metadata columns.
$ws.Range("I1").Value = "meta_title"
$ws.Range("J1").Value = "meta_keywords"
$ws.Range("K1").Value = "meta_description"

# Match the formatting of the existing header cells (font/color) on the
# newly added header cells, then wrap text across the whole header row.
$ws.Range("A1").Copy()
$ws.Range("I1:K1").PasteSpecial(-4122)
$ws.Range("A1:K1").WrapText = $true

# Add the blank (template) data row underneath, word-wrapped as well.
$ws.Range("A2:L2").WrapText = $true

$ws.Range("A6").Select() | Out-Null
